# Add a new row (PF/1.0.2) to the meta-sheet, as described in the commit
# "Add PF/1.0.2 to meta-sheet".
#
# Existing sheet layout:
#   Row 1: dev2 | sit2 | uat2 | prod
#   Row 2: PF/1.0.0 | PF/1.0.0 | PF/1.0.0 | PF/1.0.0
#
# New row 3: PF/1.0.2 | X | X | X

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "PF/1.0.2"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"
